$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.57077
$ws.Range("H2").Value = 4.71231
$ws.Range("I2").Value = 0.02582502173444737
$ws.Range("J2").Value = 0.02582502173444737
$ws.Range("O2").Value = 0.01611173663836548
$ws.Range("P2").Value = 0.01611173663836548
$ws.Range("Q2").Value = 0.10005961977
$ws.Range("R2").Value = 0.9005365779300001
$ws.Range("S2").Value = 0.0004160859488654806
$ws.Range("T2").Value = 0.0004160859488654804
$ws.Range("G3").Value = 1.57077
$ws.Range("H3").Value = 4.71231
$ws.Range("I3").Value = 0.02582502173444737
$ws.Range("J3").Value = 0.02582502173444737
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.701496333333333
$ws.Range("N3").Value = 8.104489000000001
$ws.Range("O3").Value = 0.68328279700753
$ws.Range("P3").Value = 0.68328279700753
$ws.Range("Q3").Value = 4.243429395510001
$ws.Range("R3").Value = 38.19086455959001
$ws.Range("S3").Value = 0.01764579308349346
$ws.Range("T3").Value = 0.01764579308349346
$ws.Range("G4").Value = 1.57077
$ws.Range("H4").Value = 4.71231
$ws.Range("I4").Value = 0.02582502173444737
$ws.Range("J4").Value = 0.02582502173444737
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.188504333333333
$ws.Range("N4").Value = 3.565513
$ws.Range("O4").Value = 0.3006054663541045
$ws.Range("P4").Value = 0.3006054663541044
$ws.Range("Q4").Value = 1.86686695167
$ws.Range("R4").Value = 16.80180256503
$ws.Range("S4").Value = 0.007763142702088436
$ws.Range("T4").Value = 0.007763142702088434
$ws.Range("I5").Value = 0.934831682683009
$ws.Range("J5").Value = 0.934831682683009
$ws.Range("O5").Value = 0.01611173663836548
$ws.Range("P5").Value = 0.01611173663836548
$ws.Range("S5").Value = 0.01506176187258869
$ws.Range("T5").Value = 0.01506176187258869
$ws.Range("I6").Value = 0.934831682683009
$ws.Range("J6").Value = 0.934831682683009
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.701496333333333
$ws.Range("N6").Value = 8.104489000000001
$ws.Range("O6").Value = 0.68328279700753
$ws.Range("P6").Value = 0.68328279700753
$ws.Range("Q6").Value = 153.6065403135679
$ws.Range("R6").Value = 1382.458862822111
$ws.Range("S6").Value = 0.6387544068749021
$ws.Range("T6").Value = 0.6387544068749021
$ws.Range("I7").Value = 0.934831682683009
$ws.Range("J7").Value = 0.934831682683009
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.188504333333333
$ws.Range("N7").Value = 3.565513
$ws.Range("O7").Value = 0.3006054663541045
$ws.Range("P7").Value = 0.3006054663541044
$ws.Range("Q7").Value = 67.57811829629857
$ws.Range("R7").Value = 608.2030646666871
$ws.Range("S7").Value = 0.2810155139355181
$ws.Range("T7").Value = 0.281015513935518
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.9273763333333335
$ws.Range("H8").Value = 2.782129
$ws.Range("I8").Value = 0.01524698967025436
$ws.Range("J8").Value = 0.01524698967025436
$ws.Range("O8").Value = 0.01611173663836548
$ws.Range("P8").Value = 0.01611173663836548
$ws.Range("Q8").Value = 0.05907479980966667
$ws.Range("R8").Value = 0.531673198287
$ws.Range("S8").Value = 0.0002456554820950172
$ws.Range("T8").Value = 0.0002456554820950171
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.9273763333333335
$ws.Range("H9").Value = 2.782129
$ws.Range("I9").Value = 0.01524698967025436
$ws.Range("J9").Value = 0.01524698967025436
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.701496333333333
$ws.Range("N9").Value = 8.104489000000001
$ws.Range("O9").Value = 0.68328279700753
$ws.Range("P9").Value = 0.68328279700753
$ws.Range("Q9").Value = 2.505303764120112
$ws.Range("R9").Value = 22.547733877081
$ws.Range("S9").Value = 0.01041800574783632
$ws.Range("T9").Value = 0.01041800574783632
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.9273763333333335
$ws.Range("H10").Value = 2.782129
$ws.Range("I10").Value = 0.01524698967025436
$ws.Range("J10").Value = 0.01524698967025436
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.188504333333333
$ws.Range("N10").Value = 3.565513
$ws.Range("O10").Value = 0.3006054663541045
$ws.Range("P10").Value = 0.3006054663541044
$ws.Range("Q10").Value = 1.102190790797445
$ws.Range("R10").Value = 9.919717117177001
$ws.Range("S10").Value = 0.004583328440323026
$ws.Range("T10").Value = 0.004583328440323025
$ws.Range("G11").Value = 0.7810079999999999
$ws.Range("H11").Value = 2.343024
$ws.Range("I11").Value = 0.0128405486320577
$ws.Range("J11").Value = 0.0128405486320577
$ws.Range("O11").Value = 0.01611173663836548
$ws.Range("P11").Value = 0.01611173663836548
$ws.Range("Q11").Value = 0.04975099060799999
$ws.Range("R11").Value = 0.4477589154719999
$ws.Range("S11").Value = 0.0002068835378518377
$ws.Range("T11").Value = 0.0002068835378518377
$ws.Range("G12").Value = 0.7810079999999999
$ws.Range("H12").Value = 2.343024
$ws.Range("I12").Value = 0.0128405486320577
$ws.Range("J12").Value = 0.0128405486320577
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 2.701496333333333
$ws.Range("N12").Value = 8.104489000000001
$ws.Range("O12").Value = 0.68328279700753
$ws.Range("P12").Value = 0.68328279700753
$ws.Range("Q12").Value = 2.109890248304
$ws.Range("R12").Value = 18.989012234736
$ws.Range("S12").Value = 0.008773725984423596
$ws.Range("T12").Value = 0.008773725984423596
$ws.Range("G13").Value = 0.7810079999999999
$ws.Range("H13").Value = 2.343024
$ws.Range("I13").Value = 0.0128405486320577
$ws.Range("J13").Value = 0.0128405486320577
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 1.188504333333333
$ws.Range("N13").Value = 3.565513
$ws.Range("O13").Value = 0.3006054663541045
$ws.Range("P13").Value = 0.3006054663541044
$ws.Range("Q13").Value = 0.9282313923679999
$ws.Range("R13").Value = 8.354082531311999
$ws.Range("S13").Value = 0.003859939109782262
$ws.Range("T13").Value = 0.003859939109782261
$ws.Range("G14").Value = 0.6846153333333334
$ws.Range("H14").Value = 2.053846
$ws.Range("I14").Value = 0.01125575728023152
$ws.Range("J14").Value = 0.01125575728023152
$ws.Range("O14").Value = 0.01611173663836548
$ws.Range("P14").Value = 0.01611173663836548
$ws.Range("Q14").Value = 0.04361068134866666
$ws.Range("R14").Value = 0.392496132138
$ws.Range("S14").Value = 0.0001813497969644552
$ws.Range("T14").Value = 0.0001813497969644551
$ws.Range("G15").Value = 0.6846153333333334
$ws.Range("H15").Value = 2.053846
$ws.Range("I15").Value = 0.01125575728023152
$ws.Range("J15").Value = 0.01125575728023152
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 2.701496333333333
$ws.Range("N15").Value = 8.104489000000001
$ws.Range("O15").Value = 0.68328279700753
$ws.Range("P15").Value = 0.68328279700753
$ws.Range("Q15").Value = 1.849485812743778
$ws.Range("R15").Value = 16.645372314694
$ws.Range("S15").Value = 0.007690865316874461
$ws.Range("T15").Value = 0.007690865316874461
$ws.Range("G16").Value = 0.6846153333333334
$ws.Range("H16").Value = 2.053846
$ws.Range("I16").Value = 0.01125575728023152
$ws.Range("J16").Value = 0.01125575728023152
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 1.188504333333333
$ws.Range("N16").Value = 3.565513
$ws.Range("O16").Value = 0.3006054663541045
$ws.Range("P16").Value = 0.3006054663541044
$ws.Range("Q16").Value = 0.8136682903331112
$ws.Range("R16").Value = 7.323014612998
$ws.Range("S16").Value = 0.003383542166392601
$ws.Range("T16").Value = 0.003383542166392601
